$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Rename the "Função" values: drop the underscore-separated naming in favour of
# spaces, and fix the row/value assignment so each position's B column lists the
# correct sub-function (the original file had the Lateral/Zagueiro/Primeiro
# Volante/Segundo Volante/Meia/Extremo/Atacante groups mismatched).
$ws.Range("B4").Value = "Lateral Equilibrado"
$ws.Range("B5").Value = "Lateral Ofensivo"
$ws.Range("B6").Value = "Lateral Defensivo"
$ws.Range("B7").Value = "Zagueiro Equilibrado"
$ws.Range("B8").Value = "Zagueiro Construtor"
$ws.Range("B9").Value = "Zagueiro Defensivo"
$ws.Range("B10").Value = "Primeiro Volante Equilibrado"
$ws.Range("B11").Value = "Primeiro Volante Construtor"
$ws.Range("B12").Value = "Primeiro Volante Defensivo"
$ws.Range("B13").Value = "Segundo Volante Equilibrado"
$ws.Range("B14").Value = "Segundo Volante Box to Box"
$ws.Range("B15").Value = "Segundo Volante Organizador"
$ws.Range("B16").Value = "Meia Atacante"
$ws.Range("B17").Value = "Meia Organizador"
$ws.Range("B18").Value = "Extremo Agudo"
$ws.Range("B19").Value = "Extremo Organizador"
$ws.Range("B20").Value = "Extremo Tático"
$ws.Range("B21").Value = "Atacante Referência"
$ws.Range("B22").Value = "Atacante Móvel"
$ws.Range("B23").Value = "Segundo Atacante"

# B10 picks up the same vertically-centred style already used by B11:B23.
$ws.Range("B10").VerticalAlignment = -4108

# Move the active selection from A23 to A21.
$null = $ws.Range("A21").Select()
